$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the "Main" table to include the new row (A1:H40 -> A1:H41) ---
$lo = $ws.ListObjects.Item("Main")
$lo.Resize($ws.Range("A1:H41"))

# --- Fill in the new part row (row 41) ---
$ws.Range("A41").Value = "VL53L1X ToF module"
$ws.Range("B41").Value = 4
$ws.Range("C41").Value = "Ordered"
$ws.Range("D41").Value = 17.57
$ws.Range("E41").Formula = "=PRODUCT(B41*D41)"
$ws.Range("F41").Value = 45257
$ws.Range("H41").Value = "AliExpress (SAMIORE Store)"

# --- Add the hyperlink for the LINK column, then copy the existing ---
# --- hyperlink-cell formatting so it matches the rest of the column ---
$ws.Hyperlinks.Add($ws.Range("G41"), "https://pl.aliexpress.com/item/4000074204979.html")
$ws.Range("G40").Copy()
$ws.Range("G41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the view: selection moved to H42, window scrolled up a bit ---
$ws.Range("H42").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1

$wb.Save()
